$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("A51").Value = 44854
$ws.Range("B51").Value = 318
$ws.Range("C51").Value = 119
$ws.Range("D51").Value = 107
$ws.Range("E51").Value = 259
$ws.Range("F51").Value = 212
$ws.Range("G51").Value = 5439
$ws.Range("H51").Formula = "=Data[[#This Row],[LoC]]-G50"
$ws.Range("I51").Value = 6328
$ws.Range("J51").Value = 1893
$ws.Range("K51").Value = 262
$ws.Range("L51").Value = 279
$ws.Range("M51").Value = 98
$ws.Range("N51").Value = 81
$ws.Range("O51").Value = 57
$ws.Range("P51").Formula = "=SUM(Data[[#This Row],[Shell]:[Bash]])"
$ws.Range("Q51").Formula = "=Data[[#This Row],[Total]]-P50"
$ws.Range("R51").Value = 1934
$ws.Range("S51").Value = 4078
$ws.Range("T51").Value = 65464
$ws.Range("U51").Value = 45297
$ws.Range("V51").Value = 0
$ws.Range("W51").Value = 0
$ws.Range("X51").Value = 241
$ws.Range("Y51").Formula = "=Data[[#This Row],[Open issues]]+Data[[#This Row],[Closed issues]]"
$ws.Range("Z51").Value = 0
$ws.Range("AA51").Value = 167
$ws.Range("AB51").Formula = "=Data[[#This Row],[Open pull requests]]+Data[[#This Row],[Closed pull requests]]"
$ws.Range("AC51").Value = 153
$ws.Range("AD51").Value = 164
$ws.Range("AE51").Value = 4
$ws.Range("AF51").Value = 0
$ws.Range("AG51").Value = 359
$ws.Range("AH51").Value = 913
$ws.Range("AI51").Value = 9
$ws.Range("AK51").Formula = "=SUM(Data[[#This Row],[Running]:[GH runs]])"
